$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.054.44'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.480.70'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.34'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.00'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.512'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.480.35'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.17%  '
$ws.Range("E10").Value = '  +0.87%  '
$ws.Range("E11").Value = '  -0.75%  '
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("E13").Value = '  -2.25%  '
$ws.Range("E14").Value = '  -1.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.35'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.991.36'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("E17").Value = '  -1.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.466.32'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.92'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -6.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.41'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '349.46'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.00'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.54%  '
$ws.Range("E23").Value = '  +0.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.51'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -4.08%  '
$ws.Range("E25").Value = '  -5.21%  '
$ws.Range("E26").Value = '  -2.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.24'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -8.86%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0896'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '506.94'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.46%  '
$ws.Range("E32").Value = '  -5.93%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.76'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.02%  '
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.23'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.48%  '
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '158.85'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.62%  '
$ws.Range("E37").Value = '  -8.34%  '
$ws.Range("E38").Value = '  +0.50%  '
$ws.Range("E39").Value = '  -4.47%  '
$ws.Range("E41").Value = '  -0.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.67'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.07%  '
$ws.Range("E43").Value = '  -1.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.78'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.35'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -5.04%  '
$ws.Range("E46").Value = '  -1.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '142.57'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.35%  '
$ws.Range("E48").Value = '  -4.55%  '
$ws.Range("E49").Value = '  -5.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0249'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -6.31%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0726'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.57%  '
